$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2025-08-17 08:02:26"
$ws.Range("B3").Value = 60.06000137329102
$ws.Range("C3").Value = 664.5999755859375
$ws.Range("D3").Value = 318.3999938964844
